# Translate the "formulas.xlsx" test workbook from Norwegian to English.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header row (row 1): "Kolonne N" -> "Column N" for columns B..T (N = 1..19)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $n = $i + 1
    $addr = "$($cols[$i])1"
    $ws.Range($addr).Value = "Column $n"
}

# 2. Formula / value translation in rows 2..20, columns B..T:
#    "Rad: " & ROW() & ", Kolonne: " & COLUMN()  ->  "Row: " & ROW() & ", Column: " & COLUMN()
for ($r = 2; $r -le 20; $r++) {
    $addr = "B$($r):T$($r)"
    $ws.Range($addr).Formula = '="Row: " & ROW() & ", Column: " & COLUMN()'
}

# 3. Style fix-up: move the readingOrder="0" alignment attribute from the
#    body-cell style (B2:T20) onto the two header-cell styles (B1, C1:T1).
#    Note: the readingOrder alignment attribute corresponds to xlContext (0),
#    which is the default reading-order value.
$ws.Range("B1").ReadingOrder = 0
$ws.Range("C1:T1").ReadingOrder = 0
$ws.Range("B2:T20").ReadingOrder = 1
